# Show late downpayments (agencies)
# Adds booking #15 (row 16) to the Sheet1 table and moves the selection,
# matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New booking row (row 16) ------------------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 320
$ws.Range("C16").Value = 44930
$ws.Range("D16").Value = 44932
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 44916
$ws.Range("H16").Value = "debit_card"
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 569784126

# Match the date formatting already used by the rest of the C/D/G columns.
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("G15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("C16").Value = 44930
$ws.Range("D16").Value = 44932
$ws.Range("G16").Value = 44916

# Extend the downpayment formula down into the new row (shared formula
# E4:E15 -> E4:E16).
$ws.Range("E4:E16").Formula = "= (20*B4)/100"

# The sheet was left with this selection after the edit.
$ws.Range("F20").Select()
